$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.471.44"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "3.182.10"
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.76%  "

$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.179.61"
$ws.Range("E9").Value = "  -1.18%  "

$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.514"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").Value = "3.706.18"
$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("D16").Value = "66.526.44"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "3.179.85"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.82%  "

$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("E34").Value = "  -1.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "511.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0895"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0423"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.126"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.68%  "

$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "

$ws.Range("D42").Value = "0.0₃0686"
$ws.Range("E42").Value = "  +6.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.300"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "

$ws.Range("D46").Value = "2.854.30"
$ws.Range("E46").Value = "  -5.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.27%  "

$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.21%  "
